$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in B1 from "Release Definition Name" to "Release Name"
$ws.Range("B1").Value = "Release Name"

# Update the selected/active cell to D2
$ws.Range("D2").Select()
